$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.916.42"
$ws.Cells.Item(2, 5).Value = "  +0.02%  "
$ws.Cells.Item(3, 4).Value = "1.876.99"
$ws.Cells.Item(3, 5).Value = "  -0.59%  "
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).Value = "'0.7410"
$ws.Cells.Item(5, 5).Value = "  -4.30%  "
$ws.Cells.Item(6, 4).Value = "'242.39"
$ws.Cells.Item(6, 5).Value = "  -0.21%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).Value = "'0.3153"
$ws.Cells.Item(8, 5).Value = "  +1.61%  "
$ws.Cells.Item(9, 4).Value = "'0.07209"
$ws.Cells.Item(9, 5).Value = "  +0.54%  "
$ws.Cells.Item(10, 4).Value = "'24.69"
$ws.Cells.Item(10, 5).Value = "  -3.80%  "
$ws.Cells.Item(11, 4).Value = "'0.08386"
$ws.Cells.Item(11, 5).Value = "  -2.44%  "
$ws.Cells.Item(12, 4).Value = "'0.7509"
$ws.Cells.Item(12, 5).Value = "  -1.79%  "
$ws.Cells.Item(13, 4).Value = "'5.431"
$ws.Cells.Item(13, 5).Value = "  +1.25%  "
$ws.Cells.Item(14, 4).Value = "1.876.15"
$ws.Cells.Item(14, 5).Value = "  -0.51%  "
$ws.Cells.Item(15, 4).Value = "'92.74"
$ws.Cells.Item(15, 5).Value = "  -1.13%  "
$ws.Cells.Item(16, 4).Value = "29.923.53"
$ws.Cells.Item(16, 5).Value = "  -0.04%  "
$ws.Cells.Item(17, 4).Value = "'6.079"
$ws.Cells.Item(17, 5).Value = "  -1.22%  "
$ws.Cells.Item(18, 4).Value = "'247.21"
$ws.Cells.Item(18, 5).Value = "  +1.23%  "
$ws.Cells.Item(19, 4).Value = "'13.58"
$ws.Cells.Item(19, 5).Value = "  -1.34%  "
$ws.Cells.Item(20, 4).Value = "'0.000007846"
$ws.Cells.Item(20, 5).Value = "  +0.33%  "
$ws.Cells.Item(21, 4).Value = "'0.9989"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "
$ws.Cells.Item(22, 4).Value = "2.127.23"
$ws.Cells.Item(22, 5).Value = "  -2.64%  "
$ws.Cells.Item(23, 4).Value = "'8.024"
$ws.Cells.Item(23, 5).Value = "  +0.85%  "
$ws.Cells.Item(24, 4).Value = "'1.001"
$ws.Cells.Item(24, 5).Value = "  +0.08%  "
$ws.Cells.Item(25, 4).Value = "'0.1559"
$ws.Cells.Item(25, 5).Value = "  -4.93%  "
$ws.Cells.Item(26, 4).Value = "'9.269"
$ws.Cells.Item(26, 5).Value = "  -0.96%  "
$ws.Cells.Item(28, 4).Value = "'18.62"
$ws.Cells.Item(28, 5).Value = "  -0.70%  "
$ws.Cells.Item(29, 4).Value = "'2.039"
$ws.Cells.Item(29, 5).Value = "  +0.08%  "
$ws.Cells.Item(30, 4).Value = "'1.500"
$ws.Cells.Item(30, 5).Value = "  +4.21%  "
$ws.Cells.Item(31, 4).Value = "'4.609"
$ws.Cells.Item(31, 5).Value = "  +2.29%  "
$ws.Cells.Item(32, 4).Value = "'1.535"
$ws.Cells.Item(32, 5).Value = "  +0.07%  "
$ws.Cells.Item(33, 5).Value = "  +4.15%  "
$ws.Cells.Item(34, 4).Value = "'0.05315"
$ws.Cells.Item(34, 5).Value = "  -2.25%  "
$ws.Cells.Item(35, 4).Value = "'1.239"
$ws.Cells.Item(35, 5).Value = "  -0.20%  "
$ws.Cells.Item(36, 4).Value = "'0.7544"
$ws.Cells.Item(36, 5).Value = "  +1.02%  "
$ws.Cells.Item(37, 4).Value = "'1.000"
$ws.Cells.Item(37, 5).Value = "  -0.24%  "
$ws.Cells.Item(38, 5).Value = "  -0.19%  "
$ws.Cells.Item(39, 4).Value = "'0.01961"
$ws.Cells.Item(39, 5).Value = "  -0.07%  "
$ws.Cells.Item(40, 5).Value = "  -0.90%  "
$ws.Cells.Item(41, 4).Value = "'0.4503"
$ws.Cells.Item(41, 5).Value = "  +0.84%  "
$ws.Cells.Item(42, 4).Value = "1.109.08"
$ws.Cells.Item(42, 5).Value = "  +0.11%  "
$ws.Cells.Item(43, 5).Value = "  -0.41%  "
$ws.Cells.Item(44, 4).Value = "'72.64"
$ws.Cells.Item(44, 5).Value = "  -0.60%  "
$ws.Cells.Item(45, 4).Value = "'0.8561"
$ws.Cells.Item(45, 5).Value = "  +0.79%  "
$ws.Cells.Item(46, 5).Value = "  +0.07%  "
$ws.Cells.Item(47, 4).Value = "'103.26"
$ws.Cells.Item(47, 5).Value = "  -0.39%  "
$ws.Cells.Item(48, 4).Value = "'7.629"
$ws.Cells.Item(48, 5).Value = "  +0.29%  "
$ws.Cells.Item(49, 5).Value = "  -0.73%  "
$ws.Cells.Item(50, 4).Value = "'9.523"
$ws.Cells.Item(50, 5).Value = "  -2.65%  "
$ws.Cells.Item(51, 4).Value = "2.024.72"
$ws.Cells.Item(51, 5).Value = "  -4.35%  "
